$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# Insert four new columns, positions measured against the ORIGINAL 15-column
# layout and processed left-to-right (each index already accounts for the
# columns inserted earlier in this same list):
#   col 2  -> before "Description"   => "Title Error Message"
#   col 5  -> before "SubcategoryId" => "Description Error Message"
#   col 11 -> before "EndDate"       => "StartDate Error Message"
#   col 18 -> before "State"         => "Upload File Path"
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(18).Insert()

# Populate the new cells (value-entry order chosen to reproduce the same
# shared-string table ordering as the authored workbook).
$ws.Cells.Item(1, 18).Value = "Upload File Path"
$ws.Cells.Item(2, 18).Value = "C:\Users\arjun\Desktop\dd.txt"

$ws.Cells.Item(1, 2).Value = "Title Error Message"
$ws.Cells.Item(2, 2).Value = "First character must be an alphabet character or a number."

$ws.Cells.Item(2, 5).Value = "Special characters are not allowed."
$ws.Cells.Item(1, 5).Value = "Description Error Message"

$ws.Cells.Item(1, 11).Value = "StartDate Error Message"
$ws.Cells.Item(2, 11).Value = "Start Date cannot be set to a day in the past"

# Updated StartDate / EndDate sample values.
$ws.Cells.Item(2, 10).Value = 27082019
$ws.Cells.Item(2, 12).Value = 30082019

# Column widths for the new 19-column layout (A..S), matching the authored
# widths as closely as this engine's column-width quantization allows.
$ws.Columns.Item(1).ColumnWidth = 18.3073
$ws.Columns.Item(2).ColumnWidth = 55.0221
$ws.Columns.Item(3).ColumnWidth = 35.3073
$ws.Columns.Item(4).ColumnWidth = 35.3073
$ws.Columns.Item(5).ColumnWidth = 21.5924
$ws.Columns.Item(6).ColumnWidth = 13.5924
$ws.Columns.Item(7).ColumnWidth = 8.5924
$ws.Columns.Item(8).ColumnWidth = 17.7370
$ws.Columns.Item(9).ColumnWidth = 17.7370
$ws.Columns.Item(10).ColumnWidth = 10.8776
$ws.Columns.Item(11).ColumnWidth = 39.1667
$ws.Columns.Item(12).ColumnWidth = 11.5924
$ws.Columns.Item(13).ColumnWidth = 8.7370
$ws.Columns.Item(14).ColumnWidth = 11.1667
$ws.Columns.Item(15).ColumnWidth = 15.1667
$ws.Columns.Item(16).ColumnWidth = 15.1667
$ws.Columns.Item(18).ColumnWidth = 28.5924

# Scroll position & selection to match the reviewed layout.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K2").Select()

# Page setup - portrait, paper size 9 (A4), as added during the review pass.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
